$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cells to reflect refreshed crypto price/volume data.
# Columns B/C (text) and E (percent strings) are plain text already;
# column D values are prefixed with a leading apostrophe so Excel keeps
# them as literal text instead of auto-converting to numbers (preserving
# formatting such as leading/trailing zeros, e.g. '311.60' or '0.9999').

$ws.Range("D2").Value = "'27.038.42"
$ws.Range("E2").Value = "  -1.58%  "

$ws.Range("D3").Value = "'1.828.28"
$ws.Range("E3").Value = "  -0.32%  "

$ws.Range("E4").Value = "  -0.50%  "

$ws.Range("D5").Value = "'311.60"
$ws.Range("E5").Value = "  -1.25%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.42%  "

$ws.Range("D7").Value = "'0.4309"
$ws.Range("E7").Value = "  +0.19%  "

$ws.Range("D8").Value = "'0.3666"
$ws.Range("E8").Value = "  -1.09%  "

$ws.Range("D9").Value = "'0.07278"
$ws.Range("E9").Value = "  +0.02%  "

$ws.Range("D10").Value = "'0.8448"
$ws.Range("E10").Value = "  -2.61%  "

$ws.Range("D11").Value = "'20.68"
$ws.Range("E11").Value = "  -2.63%  "

$ws.Range("D12").Value = "'1.826.96"
$ws.Range("E12").Value = "  -0.21%  "

$ws.Range("D13").Value = "'6.665"
$ws.Range("E13").Value = "  -0.57%  "

$ws.Range("D14").Value = "'0.07061"
$ws.Range("E14").Value = "  -0.39%  "

$ws.Range("D15").Value = "'5.297"

$ws.Range("D16").Value = "'89.68"
$ws.Range("E16").Value = "  +1.94%  "

$ws.Range("E17").Value = "  -0.50%  "

$ws.Range("D18").Value = "'0.000008774"
$ws.Range("E18").Value = "  -1.63%  "

$ws.Range("E19").Value = "  -0.38%  "

$ws.Range("D20").Value = "'14.92"
$ws.Range("E20").Value = "  -2.26%  "

$ws.Range("D21").Value = "'27.071.19"
$ws.Range("E21").Value = "  -1.38%  "

$ws.Range("D22").Value = "'5.142"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").Value = "'10.91"
$ws.Range("E23").Value = "  -0.05%  "

$ws.Range("D24").Value = "'2.054.04"
$ws.Range("E24").Value = "  +0.03%  "

$ws.Range("E25").Value = "  -1.19%  "

$ws.Range("D26").Value = "'151.14"
$ws.Range("E26").Value = "  -1.53%  "

$ws.Range("D27").Value = "'2.212"
$ws.Range("E27").Value = "  +2.41%  "

$ws.Range("D28").Value = "'18.28"
$ws.Range("E28").Value = "  -1.25%  "

$ws.Range("D29").Value = "'5.239"
$ws.Range("E29").Value = "  -1.04%  "

$ws.Range("D30").Value = "'117.06"
$ws.Range("E30").Value = "  -0.53%  "

$ws.Range("D31").Value = "'0.08724"
$ws.Range("E31").Value = "  -1.47%  "

$ws.Range("D32").Value = "'1.181"
$ws.Range("E32").Value = "  -2.47%  "

$ws.Range("D33").Value = "'0.7401"
$ws.Range("E33").Value = "  -3.72%  "

$ws.Range("B34").Value = "Filecoin"
$ws.Range("C34").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D34").Value = "'4.434"
$ws.Range("E34").Value = "  -1.38%  "

$ws.Range("B35").Value = "HuobiToken"
$ws.Range("C35").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D35").Value = "'2.900"
$ws.Range("E35").Value = "  -0.20%  "

$ws.Range("D36").Value = "'0.9996"
$ws.Range("E36").Value = "  -0.54%  "

$ws.Range("D37").Value = "'1.092"

$ws.Range("D38").Value = "'0.01944"
$ws.Range("E38").Value = "  -0.89%  "

$ws.Range("D39").Value = "'0.05230"
$ws.Range("E39").Value = "  -1.15%  "

$ws.Range("D40").Value = "'7.215"
$ws.Range("E40").Value = "  +0.56%  "

$ws.Range("D41").Value = "'2.866"
$ws.Range("E41").Value = "  -0.31%  "

$ws.Range("D42").Value = "'0.1701"
$ws.Range("E42").Value = "  +1.26%  "

$ws.Range("D43").Value = "'0.5143"
$ws.Range("E43").Value = "  +0.98%  "

$ws.Range("D44").Value = "'8.567"
$ws.Range("E44").Value = "  -1.42%  "

$ws.Range("D45").Value = "'10.55"
$ws.Range("E45").Value = "  -0.27%  "

$ws.Range("D46").Value = "'0.4759"
$ws.Range("E46").Value = "  +0.31%  "

$ws.Range("D47").Value = "'1.947"
$ws.Range("E47").Value = "  +6.42%  "

$ws.Range("D48").Value = "'105.79"
$ws.Range("E48").Value = "  -0.57%  "

$ws.Range("D49").Value = "'0.9990"
$ws.Range("E49").Value = "  -0.60%  "

$ws.Range("D50").Value = "'1.670"
$ws.Range("E50").Value = "  -0.15%  "

$ws.Range("D51").Value = "'0.06334"
$ws.Range("E51").Value = "  -1.52%  "
